$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp cell (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 08:35"

# --- Insert a new "Burundi" data row at row 181, shifting the former
#     Puerto Rico..Gambia block (old rows 181-189) down by one row, and
#     drop the old standalone "Burundi" row that used to sit between
#     Gambia and Granada (old row 190) since its data now lives at the
#     new row 181. ---

# Country names (column A) for rows 181-190 after the shift
$names = @(
    "Burundi",
    "Puerto Rico",
    "San Martin (Parte Francesa)",
    "Eritrea",
    "Guam",
    "Nicaragua",
    "Botsuana",
    "Antigua y Barbuda",
    "Timor Oriental",
    "Gambia"
)

# Numeric data (B,C,D,E,F,G,H) for rows 181-190 after the shift
$data = @(
    @(42, 19, 20, 21, 0, 0, 1),
    @(39, 0, 1, 36, 0, 0, 2),
    @(39, 0, 30, 6, 0, 0, 3),
    @(39, 0, 39, 0, 0, 0, 0),
    @(32, 0, 0, 31, 0, 0, 1),
    @(25, 0, 7, 10, 0, 0, 8),
    @(25, 0, 17, 7, 0, 0, 1),
    @(25, 0, 19, 3, 0, 0, 3),
    @(24, 0, 24, 0, 0, 0, 0),
    @(23, 0, 12, 10, 0, 0, 1)
)

for ($i = 0; $i -lt 10; $i++) {
    $row = 181 + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}

# --- Swap Seychelles/Groenlandia (rows 208 and 210); Montserrat (209) stays put ---
$tmp = $ws.Range("A208").Value()
$ws.Range("A208").Value = $ws.Range("A210").Value()
$ws.Range("A210").Value = $tmp

# --- Swap "Bonaire, San Eustaquio y Saba"/"San Bartolome" (rows 214 and 216);
#     Sahara Occidental (215) stays put ---
$tmp2 = $ws.Range("A214").Value()
$ws.Range("A214").Value = $ws.Range("A216").Value()
$ws.Range("A216").Value = $tmp2
